# ============================================================
# Budget.xlsx edit: restructure the second (per-item) table,
# add a "Total cost" column with per-row + aggregate formulas,
# add a cost-projection block, reposition the picture, and
# tidy up view/column-width cosmetics.
# ============================================================

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------
# 1) Remove the rows that disappear entirely from table 2
#    (delete bottom-up so row numbers stay stable as we go)
#    Row 25 = Lighting system (2 mini headlights)
#    Row 23 = Buzzer
#    Row 22 = Rechargeable Battery system
#    Row 19 = Raspberry Pi Zero SD Card (8GB)
#    Row 18 = Raspberry Pi Zero Power supply
# ------------------------------------------------------------
$ws.Rows(25).EntireRow.Delete()
$ws.Rows(23).EntireRow.Delete()
$ws.Rows(22).EntireRow.Delete()
$ws.Rows(19).EntireRow.Delete()
$ws.Rows(18).EntireRow.Delete()

# After the deletions above, table 2 rows are now:
#  15 header / 16 Raspberry Pi Zero / 17 Raspberry Pi Zero Camera /
#  18 Car platform / 19 Car Rechargeable Battery / 20 Odometric Sensors /
#  21 STM32 / 22 Smartphone (RAM=?, Storage=?, SO=?)

# ------------------------------------------------------------
# 2) Re-label / re-value the remaining rows to match the new data
# ------------------------------------------------------------
$ws.Range("A16").Value = "Remote Vision susbsystem"
$ws.Range("B16").Value = 5

$ws.Range("A17").Value = "Camera"
$ws.Range("B17").Value = 20

$ws.Range("A21").Value = "Navigation Subsystem"

$ws.Range("A22").Value = "Smartphone"
$ws.Range("C22").Value = $null

$ws.Range("B20").Value = 7.95
$ws.Range("E20").Value = 9

# ------------------------------------------------------------
# 3) Header row: rename Value-> Unit cost, add Total cost header
# ------------------------------------------------------------
$ws.Range("B15").Value = "Unit cost (€)"
$ws.Range("G15").Value = "Total cost (€)"
$ws.Range("A15:F15").Copy()
$ws.Range("G15").PasteSpecial(-4122)

# ------------------------------------------------------------
# 4) New "Total cost" column (G) for each item row: =Unit*Qty
#    Copy the tan-fill format used by the rest of that row block
#    into G, then write the formulas.
# ------------------------------------------------------------
$ws.Range("G16").Formula = "=B16*E16"
$ws.Range("G17").Formula = "=B17*E17"
$ws.Range("G18").Formula = "=B18*E18"
$ws.Range("G19").Formula = "=B19*E19"
$ws.Range("G20").Formula = "=B20*E20"
$ws.Range("G21").Formula = "=B21*E21"
$ws.Range("G22").Formula = "=B22*E22"

$ws.Range("G16").Copy()
$ws.Range("G17:G22").PasteSpecial(-4122)

# Total (sum) row right under the table
$ws.Range("G16").Copy()
$ws.Range("G23").PasteSpecial(-4122)
$ws.Range("G23").Font.Bold = $true
$ws.Range("G23").Formula = "=SUM(G16:G22)"

# ------------------------------------------------------------
# 5) Cost projection block (rows 25-28)
# ------------------------------------------------------------
$ws.Range("D25").Value = 10000
$ws.Range("G25").Formula = "=G23*D25"

$ws.Range("D26").Value = 250
$ws.Range("E26").Value = 7
$ws.Range("F26").Value = 10
$ws.Range("G26").Formula = "=D26*E26*F26"

$ws.Range("G27").Formula = "=SUM(G25+G26)"

$ws.Range("G28").Formula = "=G27/10000"

# ------------------------------------------------------------
# 6) Fix the hyperlink: it now lives on the Car platform row (C18)
# ------------------------------------------------------------
$ws.Range("C20").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C18"), "https://www.botnroll.com/pt/bases/2151--kit-iniciacao-para-robo-mecanica-.html")

# ------------------------------------------------------------
# 7) Column A width + selection/view tidy-up
# ------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 21.8
$ws.Range("H12").Select()

# ------------------------------------------------------------
# 8) Reposition / resize the picture next to the new table
# ------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Left = 754.315748031496
$shp.Top = 43.67251968503937
$shp.Width = 218.16133858267716
$shp.Height = 166.68425196850393
